$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 49, 48, 41, 35 (in descending order so earlier deletions
# don't shift the row numbers of rows still to be deleted)
$ws.Rows.Item(49).Delete()
$ws.Rows.Item(48).Delete()
$ws.Rows.Item(41).Delete()
$ws.Rows.Item(35).Delete()
